# pptx: Include image title in description
#
# Previously, when a Markdown image used a title, e.g.
#   ![alt text](link "title")
# the title was dropped when writing to pptx. The picture's description
# (PowerPoint's "Alt Text" / the OOXML <p:cNvPr descr="..."/> attribute)
# only contained the link. Now the title is included too, by prefixing
# the description with "fig:  " (the marker also used to flag that the
# description accompanies a figure reference).
#
# Walk every slide/shape, and for each picture whose description is the
# plain "lalune.jpg" link, rewrite it to include the "fig:" title prefix.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.Type -eq 13) {
            if ($sh.AlternativeText -eq "lalune.jpg") {
                $sh.AlternativeText = "fig:  lalune.jpg"
            }
        }
    }
}
